$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell P1
$ws.Range("P1").Value = "PrimaryImageURL"

# Autofit the new column so its width matches the "bestFit" width Excel
# computes for the new header text
$ws.Columns.Item(16).EntireColumn.AutoFit() | Out-Null

# Update the active selection to match the recorded cursor position
$ws.Range("I11").Select() | Out-Null
